# Release v0.1.0-beta: Fix validation errors and update canonical URL
$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")

# Plain text updates (safe to assign directly - Excel won't reinterpret the type)
$wsMeta.Range("B3").Value = "0.1.0"
$wsMeta.Range("B6").Value = "draft"
$wsMeta.Range("B8").Value = "2025-12-26T14:13:58+00:00"
$wsMeta.Range("B11").Value = "Value set for patient housing status"

# "false" looks like a Boolean literal to Excel's auto-detection, which would
# store it as a boolean cell (t="b") instead of text. Build it as a text
# formula result in a scratch cell, then paste the *value* back in - this
# keeps the cell's type as a shared string without leaving a stray
# quote-prefixed style behind.
$wsMeta.Range("Z99").Formula = "=""false"""
$wsMeta.Range("Z99").Copy()
$wsMeta.Cells.Item(7, 2).PasteSpecial(-4163)
$wsMeta.Range("Z99").Clear()

$wsInclude = $wb.Worksheets.Item("Include #0")

# "266935003" / "224224003" are purely numeric strings, which Excel's
# auto-detection would store as numbers. Same text-formula trick as above to
# keep them as text (matching how the existing SNOMED codes are stored).
$wsInclude.Range("Z99").Formula = "=""266935003"""
$wsInclude.Range("Z99").Copy()
$wsInclude.Cells.Item(2, 1).PasteSpecial(-4163)
$wsInclude.Range("Z99").Clear()

$wsInclude.Range("B2").Value = "Housing lack"

$wsInclude.Range("Z99").Formula = "=""224224003"""
$wsInclude.Range("Z99").Copy()
$wsInclude.Cells.Item(3, 1).PasteSpecial(-4163)
$wsInclude.Range("Z99").Clear()

$wsInclude.Range("B3").Value = "Lives in staffed home"
